$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "contoh"
$ws.Range("B2").Value = "contoh"
$ws.Range("I2").Value = "contoh@gmail.com"
$ws.Range("H2").Value = "+62 contoh"

$ws.Range("H3").Select()
